# edit.ps1 - applies the "added id to NewPizza and updated powerpoint" change:
#   1. Updates the cached "today" date field (datetimeFigureOut) from
#      5/12/2020 -> 5/26/2020 on the Slide Master and every Slide Layout.
#   2. Rewrites the "Submit Order saves..." callout textbox on slide 2 so it
#      reads "...the new row added" (dropping the Wingdings smiley + trailing
#      "(after showing "Pizza time" )") and lets the shape re-autofit to its
#      new (shorter) height.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder text: 5/12/2020 -> 5/26/2020
# ---------------------------------------------------------------------------
$newDate = "5/26/2020"

function Set-DatePlaceholderText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide Master
Set-DatePlaceholderText $p.SlideMaster.Shapes

# Every Slide Layout (CustomLayouts) hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. "Submit Order saves..." textbox on slide 2 (shape "TextBox 15")
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$box = $null
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $candidate = $slide2.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 15") {
        $box = $candidate
    }
}
if ($box -eq $null) {
    $box = $slide2.Shapes.Item(7)
}
$tr = $box.TextFrame.TextRange

# Locate "row added" in the existing sentence so the char offsets below are
# derived from the live text instead of being hard-coded twice.
$rowAddedHit = $tr.Find("row added", 1)
$rowAddedStart = $rowAddedHit.Start
$rowAddedLen = $rowAddedHit.Length
$afterRowAddedStart = $rowAddedStart + $rowAddedLen

# Drop the trailing Wingdings smiley run (the last two characters: the
# symbol char + the closing parenthesis).
$len = $tr.Length
$symRun = $tr.Characters($len - 1, 2)
$symRun.Delete()

# Drop the now-trailing ' (after showing "Pizza time" ' text, which was
# part of the first run, keeping up through "...the new row added".
$afterRowAdded = $tr.Characters($afterRowAddedStart, $tr.Length - $afterRowAddedStart + 1)
$afterRowAdded.Delete()

# Re-create "row added" as its own run so it lines up with the new text
# ("...the new " / "row added").
$rowAdded = $tr.Characters($rowAddedStart, $rowAddedLen)
$rowAdded.Delete()
$tr.InsertAfter("row added") | Out-Null
